$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update cryptocurrency price/volume data (GitHub Actions scheduled refresh)

$ws.Range("D2").Value = '49.661.15'
$ws.Range("E2").Value = '  -0.80%  '
$ws.Range("D3").Value = '2.647.81'
$ws.Range("E3").Value = '  -0.05%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.00'
$ws.Range("E4").Value = '  +0.06%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '112.68'
$ws.Range("E5").Value = '  -1.25%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '327.87'
$ws.Range("E6").Value = '  +0.41%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.523'
$ws.Range("E7").Value = '  -1.18%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '1.00'
$ws.Range("E8").Value = '  +0.03%  '
$ws.Range("E9").Value = '  -1.24%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '39.74'
$ws.Range("E10").Value = '  -3.14%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '19.98'
$ws.Range("E11").Value = '  -1.07%  '
$ws.Range("E12").Value = '  -0.77%  '
$ws.Range("E13").Value = '  +2.23%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '7.60'
$ws.Range("E14").Value = '  +2.69%  '
$ws.Range("D15").Value = '3.065.77'
$ws.Range("E15").Value = '  +0.08%  '
$ws.Range("D16").Value = '2.645.52'
$ws.Range("E16").Value = '  -0.01%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.859'
$ws.Range("E17").Value = '  -1.63%  '
$ws.Range("D18").Value = '49.662.48'
$ws.Range("E18").Value = '  -0.61%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '13.33'
$ws.Range("E19").Value = '  +0.60%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '6.69'
$ws.Range("E20").Value = '  -1.55%  '
$ws.Range("E21").Value = '  -0.40%  '
$ws.Range("D22").Value = '0.0₃0950'
$ws.Range("E22").Value = '  -0.86%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '268.88'
$ws.Range("E23").Value = '  -2.74%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '69.10'
$ws.Range("E24").Value = '  -4.20%  '
$ws.Range("E25").Value = '  -0.52%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '26.15'
$ws.Range("E26").Value = '  -2.64%  '
$ws.Range("E27").Value = '  +0.00%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '10.20'
$ws.Range("E28").Value = '  +1.43%  '
$ws.Range("E29").Value = '  -0.82%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '0.139'
$ws.Range("E30").Value = '  -1.58%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '34.87'
$ws.Range("E31").Value = '  -3.93%  '
$ws.Range("E32").Value = '  -1.42%  '
$ws.Range("E33").Value = '  +0.22%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.0821'
$ws.Range("E34").Value = '  +1.00%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '19.21'
$ws.Range("E35").Value = '  -1.68%  '
$ws.Range("E36").Value = '  -0.13%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '4.95'
$ws.Range("E37").Value = '  -2.47%  '
$ws.Range("E38").Value = '  -1.82%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '3.15'
$ws.Range("E39").Value = '  +1.28%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '128.94'
$ws.Range("E40").Value = '  +4.48%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '23.49'
$ws.Range("E41").Value = '  +6.54%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.0344'
$ws.Range("E42").Value = '  +8.51%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '2.28'
$ws.Range("E44").Value = '  -0.74%  '
$ws.Range("B45").Value = 'NEARProtocol'
$ws.Range("C45").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '3.33'
$ws.Range("E45").Value = '  -0.20%  '
$ws.Range("B46").Value = 'Maker'
$ws.Range("C46").Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range("D46").Value = '2.065.31'
$ws.Range("E46").Value = '  -1.00%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '2.11'
$ws.Range("E47").Value = '  +5.96%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '2.23'
$ws.Range("E48").Value = '  -3.12%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '8.95'
$ws.Range("E49").Value = '  -2.27%  '
$ws.Range("E50").Value = '  -2.60%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '58.92'
$ws.Range("E51").Value = '  -1.63%  '
